$wb = $excel.ActiveWorkbook

# Sheet "current": add rows 17-19 (ID 16-18), same B/C/D as previous rows
$ws = $wb.Worksheets.Item("current")
$startRow = 17
$startId = 16
for ($i = 0; $i -lt 3; $i++) {
    $r = $startRow + $i
    $id = $startId + $i
    $ws.Cells.Item($r, 1).Value2 = $id
    $ws.Cells.Item($r, 2).Value2 = "(Sat, 09 Jul 2022 09:58:02 GMT)"
    $ws.Cells.Item($r, 3).Value2 = $ws.Cells.Item(2, 3).Value2
    $ws.Cells.Item($r, 4).Value2 = "current"
}

# Sheet "forecast": add rows 10-12 (ID 9-11)
$ws = $wb.Worksheets.Item("forecast")
$startRow = 10
$startId = 9
for ($i = 0; $i -lt 3; $i++) {
    $r = $startRow + $i
    $id = $startId + $i
    $ws.Cells.Item($r, 1).Value2 = $id
    $ws.Cells.Item($r, 2).Value2 = "(Sat, 09 Jul 2022 09:58:02 GMT)"
    $ws.Cells.Item($r, 3).Value2 = $ws.Cells.Item(2, 3).Value2
    $ws.Cells.Item($r, 4).Value2 = "forecast"
}

# Sheet "football": add rows 10-12 (ID 9-11)
$ws = $wb.Worksheets.Item("football")
$startRow = 10
$startId = 9
for ($i = 0; $i -lt 3; $i++) {
    $r = $startRow + $i
    $id = $startId + $i
    $ws.Cells.Item($r, 1).Value2 = $id
    $ws.Cells.Item($r, 2).Value2 = "(Sat, 09 Jul 2022 09:58:02 GMT)"
    $ws.Cells.Item($r, 3).Value2 = $ws.Cells.Item(2, 3).Value2
    $ws.Cells.Item($r, 4).Value2 = "football"
}

# Sheet "astronomy": add rows 10-12 (ID 9-11)
$ws = $wb.Worksheets.Item("astronomy")
$startRow = 10
$startId = 9
for ($i = 0; $i -lt 3; $i++) {
    $r = $startRow + $i
    $id = $startId + $i
    $ws.Cells.Item($r, 1).Value2 = $id
    $ws.Cells.Item($r, 2).Value2 = "(Sat, 09 Jul 2022 09:58:02 GMT)"
    $ws.Cells.Item($r, 3).Value2 = $ws.Cells.Item(2, 3).Value2
    $ws.Cells.Item($r, 4).Value2 = "astronomy"
}

# Sheet "timezone": add rows 10-12 (ID 9-11)
$ws = $wb.Worksheets.Item("timezone")
$startRow = 10
$startId = 9
for ($i = 0; $i -lt 3; $i++) {
    $r = $startRow + $i
    $id = $startId + $i
    $ws.Cells.Item($r, 1).Value2 = $id
    $ws.Cells.Item($r, 2).Value2 = "(Sat, 09 Jul 2022 09:58:02 GMT)"
    $ws.Cells.Item($r, 3).Value2 = $ws.Cells.Item(2, 3).Value2
    $ws.Cells.Item($r, 4).Value2 = "timezone"
}
